# Cat, Parrot and Seed - add answer 5b after the existing "Finally, ..." line
# in the "Choose a solution and develop a plan to implement it" paragraph.
#
# Word represents a manual line break (Shift+Enter, <w:br/>) in Find/Replace
# and range text as Chr(11) (vertical-tab). We build the addition (a blank
# line followed by the new "b." sentence) out of that character so Word
# serializes real <w:br/> elements instead of literal text.

$d = $word.ActiveDocument

$lineBreak = [char]11

$existingText = "Finally, take the parrot to side B to be with the cat and the bag of seed. "
$newSentence = "b.  As long as the parrot is not left alone with either the cat or the bag of seed, the solution works out.  By moving the parrot back and forth, this solution is achieved.  "

$replacementText = $existingText + $lineBreak + $lineBreak + $newSentence

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute($existingText, $false, $false, $false, $false, $false, $true, 1, $false, $replacementText, 2)

if (-not $found) {
    throw "Could not find the target sentence to update."
}
